$d = $word.ActiveDocument

# --- Paragraph 1 originally holds two runs:
#       1) "**ID__AFFARS_5341_topic_4__ID**"
#       2) " "            (a lone trailing space)
#     The target edit:
#       - renames the placeholder id to **ID__AFFARS_SUBPART_5341_2__ID**
#       - removes the now unused trailing-space run
#       - adds a (space-only, line-less) paragraph border
#       - bumps the left indent from 120 -> 225 twips

# 1) Update the placeholder text in-place (keeps existing run formatting).
$d.Content.Find.Execute(
    "**ID__AFFARS_5341_topic_4__ID**", $true, $false, $false, $false, $false,
    $true, 1, $false, "**ID__AFFARS_SUBPART_5341_2__ID**", 2)

# 2) Drop the trailing-space run that followed the placeholder text.
$p1 = $d.Paragraphs.Item(1)
$endPos = $p1.Range.End
$spaceRange = $d.Range($endPos - 2, $endPos - 1)
$spaceRange.Delete()

# 3) Add the paragraph border + new indent. Word's Borders COM object can only
#    emit "real" (val/sz) borders, never the bare w:space-only form used here,
#    so instead insert a tiny OOXML fragment carrying the desired <w:pPr> at
#    the paragraph mark, then merge the resulting (now extra) paragraph back
#    into paragraph 1 by deleting the paragraph break between them. Because a
#    paragraph-merge keeps the *second* paragraph's properties, the inserted
#    <w:pPr> (border + spacing + indent + alignment) ends up governing the
#    merged paragraph 1, while its run(s) are left completely untouched.
$p1 = $d.Paragraphs.Item(1)
$collapsed = $d.Range($p1.Range.End - 1, $p1.Range.End - 1)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pBdr>
<w:top w:space="5"/>
<w:left w:space="5"/>
<w:bottom w:space="5"/>
<w:right w:space="5"/>
</w:pBdr>
<w:spacing w:after="0"/>
<w:ind w:left="225"/>
<w:jc w:val="left"/>
</w:pPr>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$collapsed.InsertXML($xml)

$p1 = $d.Paragraphs.Item(1)
$mergeRange = $d.Range($p1.Range.End - 1, $p1.Range.End)
$mergeRange.Delete()
